# Appends the new "intro" dialog paragraphs after the existing final
# paragraph ("You will certainly have no problem dealing with the
# remaining blobs!"), per the commit "more intro stuff, implement intro".
#
# Add-Paragraph inserts one brand-new paragraph at the end of the
# document (InsertParagraphAfter, mirroring how Word itself creates a
# paragraph mark) and then types each element of -Runs into it in turn.
# Re-fetching $d.Paragraphs.Last.Range after every mutation keeps the
# Range in sync with the growing story, since a stale Range object does
# not auto-advance here.
#
# Typing several consecutive InsertAfter runs with identical formatting
# normally coalesces them into a single <w:r> on save. Wrapping the typing
# in a throwaway TrackRevisions session (then immediately accepting it)
# keeps each InsertAfter call as its own <w:r>, matching how the source
# document represents multi-run paragraphs, without leaving any <w:ins>/
# revision markup behind in the saved file.

function Add-Paragraph {
    param($d, [string[]]$Runs)

    $r = $d.Paragraphs.Last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $d.TrackRevisions = $true
    foreach ($t in $Runs) {
        $r = $d.Paragraphs.Last.Range
        $r.Collapse(0)
        $r.InsertAfter($t)
    }
    $d.TrackRevisions = $false
    $d.AcceptAllRevisions() | Out-Null
}

$d = $word.ActiveDocument

# Blank paragraph separating the old ending from the new intro block.
Add-Paragraph $d @()

Add-Paragraph $d @(
    "Multiple space blobs have pierced through ",
    "our dimension",
    "!"
)

Add-Paragraph $d @("Emergency protocol initiated.")

Add-Paragraph $d @("We must banish them immediately before they fall down to Earth!")

Add-Paragraph $d @(
    "With our latest advancements in ",
    "blobology",
    ", we will be deploying Attack Blobs."
)

Add-Paragraph $d @(
    "These blobs must be made with the power of multiplication, and who better to ",
    "do it",
    " than you!"
)

Add-Paragraph $d @(
    "Our i",
    "ntrepid hero, ",
    "go forth",
    ", and ",
    "use",
    " ",
    "your ",
    "math",
    "ematical might to banish these invading blobs back to their dimension!"
)

Add-Paragraph $d @(" ")
